$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 15801.75
$ws.Range("I11").Value = 15801.75
$ws.Range("K11").Value = 15801.75
$ws.Range("M11").Value = -15661.75
$ws.Range("H17").Value = 189916.6
$ws.Range("J17").Value = 195472.97
$ws.Range("L17").Value = 586418.91
$ws.Range("N17").Value = -586754.91
$ws.Range("H28").Value = 582.2174
$ws.Range("I28").Value = 600.3158
$ws.Range("K28").Value = 600.3158
$ws.Range("M28").Value = -115.3158
$ws.Range("H51").Value = 45458564
$ws.Range("J51").Value = 71431310
$ws.Range("L51").Value = 71431310
$ws.Range("N51").Value = -71432278
$ws.Range("H62").Value = 4992.25
$ws.Range("I62").Value = 4992.25
$ws.Range("K62").Value = 4992.25
$ws.Range("M62").Value = -4368.25
$ws.Range("H65").Value = 4992.25
$ws.Range("I65").Value = 4992.25
$ws.Range("K65").Value = 24961.25
$ws.Range("M65").Value = -21841.25
$ws.Range("H92").Value = 2604760
$ws.Range("I92").Value = 1042252.3
$ws.Range("K92").Value = 1042252.3
$ws.Range("M92").Value = -1041004.3
$ws.Range("H100").Value = 1666
$ws.Range("I100").Value = 1666
$ws.Range("K100").Value = 1666
$ws.Range("M100").Value = -1125
$ws.Range("H112").Value = 3106.9119
$ws.Range("I112").Value = 2250
$ws.Range("K112").Value = 6750
$ws.Range("M112").Value = -5642
$ws.Range("H131").Value = 6342
$ws.Range("I131").Value = 6342
$ws.Range("K131").Value = 19026
$ws.Range("M131").Value = -13986
$ws.Range("H135").Value = 735.45
$ws.Range("I135").Value = 735.45
$ws.Range("K135").Value = 6619.05
$ws.Range("M135").Value = -4084.05
$ws.Range("H137").Value = 2047065.4
$ws.Range("I137").Value = 3867
$ws.Range("J137").Value = 7705153
$ws.Range("K137").Value = 11601
$ws.Range("L137").Value = 23115459
$ws.Range("M137").Value = -9051
$ws.Range("N137").Value = -23120559
$ws.Range("H138").Value = 6439.485
$ws.Range("I138").Value = 10967.286
$ws.Range("J138").Value = 3103.2104
$ws.Range("K138").Value = 32901.858
$ws.Range("L138").Value = 9309.6312
$ws.Range("M138").Value = -27761.858
$ws.Range("N138").Value = -19589.6312
$ws.Range("H141").Value = 1877.4736
$ws.Range("I141").Value = 1704
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 5112
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 68
$ws.Range("N141").Value = -25360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1584
$ws.Range("I2").Value = 1498
$ws.Range("K2").Value = 1498
$ws.Range("M2").Value = -1385
$ws.Range("H32").Value = 148266.08
$ws.Range("I32").Value = 164248.6
$ws.Range("J32").Value = 8989.857
$ws.Range("K32").Value = 164248.6
$ws.Range("L32").Value = 8989.857
$ws.Range("N32").Value = -9563.857
$ws.Range("M32").Value = -163961.6
$ws.Range("H45").Value = 73563.57000000001
$ws.Range("I45").Value = 92828.63
$ws.Range("K45").Value = 92828.63
$ws.Range("M45").Value = -92451.63
$ws.Range("H61").Value = 1002313.06
$ws.Range("J61").Value = 3525295
$ws.Range("L61").Value = 3525295
$ws.Range("N61").Value = -3525719
$ws.Range("H74").Value = 390304.84
$ws.Range("I74").Value = 2815.5405
$ws.Range("J74").Value = 1013657.2
$ws.Range("K74").Value = 2815.5405
$ws.Range("L74").Value = 1013657.2
$ws.Range("M74").Value = -1941.5405
$ws.Range("N74").Value = -1015405.2
$ws.Range("H77").Value = 390304.84
$ws.Range("I77").Value = 2815.5405
$ws.Range("J77").Value = 1013657.2
$ws.Range("K77").Value = 14077.7025
$ws.Range("L77").Value = 5068286
$ws.Range("M77").Value = -9709.702499999999
$ws.Range("N77").Value = -5077022
$ws.Range("H110").Value = 961.7143
$ws.Range("I110").Value = 962.9474
$ws.Range("K110").Value = 962.9474
$ws.Range("M110").Value = 1082.0526
$ws.Range("H116").Value = 1584
$ws.Range("I116").Value = 1498
$ws.Range("K116").Value = 1498
$ws.Range("M116").Value = 796
$ws.Range("H121").Value = 79998.5
$ws.Range("J121").Value = 79998.5
$ws.Range("L121").Value = 79998.5
$ws.Range("N121").Value = -83492.5
$ws.Range("H122").Value = 2114.875
$ws.Range("I122").Value = 1702.7142
$ws.Range("K122").Value = 5108.142599999999
$ws.Range("M122").Value = -2658.142599999999
$ws.Range("H132").Value = 1701.2106
$ws.Range("I132").Value = 1151.7291
$ws.Range("J132").Value = 4631.778
$ws.Range("K132").Value = 3455.1873
$ws.Range("L132").Value = 13895.334
$ws.Range("M132").Value = -925.1873000000001
$ws.Range("N132").Value = -18955.334
$ws.Range("H136").Value = 1002313.06
$ws.Range("J136").Value = 3525295
$ws.Range("L136").Value = 10575885
$ws.Range("N136").Value = -10580985

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1584
$ws.Range("I3").Value = 1498
$ws.Range("K3").Value = 1498
$ws.Range("M3").Value = -1384
$ws.Range("H22").Value = 2192.7144
$ws.Range("I22").Value = 1708.1666
$ws.Range("K22").Value = 1708.1666
$ws.Range("M22").Value = -1535.1666
$ws.Range("H105").Value = 10193.193
$ws.Range("I105").Value = 7339.143
$ws.Range("K105").Value = 7339.143
$ws.Range("M105").Value = -5592.143
$ws.Range("H134").Value = 20456168
$ws.Range("I134").Value = 1268.8928
$ws.Range("K134").Value = 3806.6784
$ws.Range("M134").Value = -1271.6784

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 946.2857
$ws.Range("J22").Value = 745
$ws.Range("L22").Value = 745
$ws.Range("N22").Value = -1445
$ws.Range("H56").Value = 25000
$ws.Range("I56").Value = 25000
$ws.Range("K56").Value = 25000
$ws.Range("M56").Value = -24155
$ws.Range("H58").Value = 4000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 4000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -4406
$ws.Range("H99").Value = 5001525
$ws.Range("I99").Value = 5001525
$ws.Range("K99").Value = 5001525
$ws.Range("M99").Value = -5000027
$ws.Range("H126").Value = 5001525
$ws.Range("I126").Value = 5001525
$ws.Range("K126").Value = 15004575
$ws.Range("M126").Value = -15002105
$ws.Range("H132").Value = 2302.077
$ws.Range("I132").Value = 1676.6
$ws.Range("J132").Value = 2693
$ws.Range("K132").Value = 5029.799999999999
$ws.Range("L132").Value = 8079
$ws.Range("M132").Value = -2499.799999999999
$ws.Range("N132").Value = -13139
$ws.Range("H134").Value = 2743.875
$ws.Range("I134").Value = 2454.818
$ws.Range("J134").Value = 3379.8
$ws.Range("K134").Value = 7364.454000000001
$ws.Range("L134").Value = 10139.4
$ws.Range("M134").Value = -4829.454000000001
$ws.Range("N134").Value = -15209.4
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 12000
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -17100

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2749.5
$ws.Range("J80").Value = 2999
$ws.Range("L80").Value = 8997
$ws.Range("N80").Value = -10869
$ws.Range("H83").Value = 2749.5
$ws.Range("J83").Value = 2999
$ws.Range("L83").Value = 26991
$ws.Range("N83").Value = -36351
$ws.Range("H122").Value = 11112414
$ws.Range("I122").Value = 13333799
$ws.Range("J122").Value = 5494
$ws.Range("K122").Value = 120004191
$ws.Range("L122").Value = 49446
$ws.Range("M122").Value = -120001741
$ws.Range("N122").Value = -54346
$ws.Range("H139").Value = 3390.0557
$ws.Range("I139").Value = 2214.2144
$ws.Range("K139").Value = 6642.6432
$ws.Range("M139").Value = -1502.6432

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2707.6667
$ws.Range("I126").Value = 1830.25
$ws.Range("J126").Value = 4462.5
$ws.Range("K126").Value = 5490.75
$ws.Range("L126").Value = 13387.5
$ws.Range("M126").Value = -3020.75
$ws.Range("N126").Value = -18327.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 679.65625
$ws.Range("I55").Value = 487
$ws.Range("K55").Value = 487
$ws.Range("M55").Value = -314
$ws.Range("H122").Value = 3732.5217
$ws.Range("I122").Value = 3132.2354
$ws.Range("K122").Value = 9396.706200000001
$ws.Range("M122").Value = -6946.706200000001
$ws.Range("H136").Value = 94261.17999999999
$ws.Range("I136").Value = 254700.75
$ws.Range("J136").Value = 2581.4285
$ws.Range("K136").Value = 764102.25
$ws.Range("L136").Value = 7744.2855
$ws.Range("M136").Value = -761552.25
$ws.Range("N136").Value = -12844.2855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1906752.2
$ws.Range("I107").Value = 1595.1
$ws.Range("K107").Value = 4785.299999999999
$ws.Range("M107").Value = -2865.299999999999
$ws.Range("H113").Value = 823.8461
$ws.Range("I113").Value = 710.7143
$ws.Range("K113").Value = 2132.1429
$ws.Range("M113").Value = 37.85710000000017
$ws.Range("H122").Value = 2020.8572
$ws.Range("I122").Value = 1954.3636
$ws.Range("J122").Value = 2264.6667
$ws.Range("K122").Value = 5863.0908
$ws.Range("L122").Value = 6794.000100000001
$ws.Range("M122").Value = -3413.0908
$ws.Range("N122").Value = -11694.0001
$ws.Range("H126").Value = 13159402
$ws.Range("I126").Value = 14707097
$ws.Range("K126").Value = 44121291
$ws.Range("M126").Value = -44118821
$ws.Range("H136").Value = 1043.25
$ws.Range("I136").Value = 889.4286
$ws.Range("J136").Value = 1258.6
$ws.Range("K136").Value = 2668.2858
$ws.Range("L136").Value = 3775.8
$ws.Range("M136").Value = -118.2857999999997
$ws.Range("N136").Value = -8875.799999999999
